# Applies the "lamda new init guesses" edit to sheets lamda0.05, lamda1, lamda0.01, lamda0.5, lamda0.15
# For each of those sheets:
#   - rows 2-4 (gain, proximity, lamda) keep their labels but get updated Mean/95%LL/95%UL values
#   - two new rows are inserted at row 5 for "pi" and "Keq" (now fully populated with Mean/LL/UL/True)
#   - the previously-last "pi" row (old row 14) is removed, since its content now lives in the new row 5
#   - the rows in between (SNR..p(specific)) keep their relative order, shifted down by two rows

$wb = $excel.ActiveWorkbook

# ---- sheet: lamda0.05 ----
$ws = $wb.Worksheets.Item("lamda0.05")

# Insert two blank rows before row 5 (will become "pi" and "Keq")
$ws.Range("A5:E6").Insert()

# Remove the old trailing "pi" row, now shifted down to row 16
$ws.Rows.Item(16).Delete()

# Give the two new label cells (A5:A6) the same bold/bordered format as the
# other row labels in column A (copy format from A4, the "lamda" label cell)
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

$ws.Range("A2").Value = "gain"
$ws.Range("B2").Value = 6.996687583894611
$ws.Range("C2").Value = 6.976823627615556
$ws.Range("D2").Value = 7.017785375982857
$ws.Range("E2").Value = 7

$ws.Range("A3").Value = "proximity"
$ws.Range("B3").Value = 0.204292133285748
$ws.Range("C3").Value = 0.1936113303745512
$ws.Range("D3").Value = 0.2138933491786145
$ws.Range("E3").Value = 0.2

$ws.Range("A4").Value = "lamda"
$ws.Range("B4").Value = 0.05422331215465174
$ws.Range("C4").Value = 0.0478392138545554
$ws.Range("D4").Value = 0.0607639377558007
$ws.Range("E4").Value = 0.05

$ws.Range("A5").Value = "pi"
$ws.Range("B5").Value = 0.1731204375081504
$ws.Range("C5").Value = 0.1587147316212725
$ws.Range("D5").Value = 0.1893954936444539
$ws.Range("E5").Value = 0.15

$ws.Range("A6").Value = "Keq"
$ws.Range("B6").Value = 0.2104955595729849
$ws.Range("C6").Value = 0.1886574486983213
$ws.Range("D6").Value = 0.2336472656222764
$ws.Range("E6").Value = ""

$ws.Range("A7").Value = "SNR"
$ws.Range("B7").Value = 3.761977261134855
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 3.75889905491011

$ws.Range("A8").Value = "MCC"
$ws.Range("B8").Value = 0.9914675955572496
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""

$ws.Range("A9").Value = "Recall"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""

$ws.Range("A10").Value = "Precision"
$ws.Range("B10").Value = 0.9858490566037735
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""

$ws.Range("A11").Value = "TN"
$ws.Range("B11").Value = 2076
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""

$ws.Range("A12").Value = "FP"
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""

$ws.Range("A13").Value = "FN"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""

$ws.Range("A14").Value = "TP"
$ws.Range("B14").Value = 418
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""

$ws.Range("A15").Value = "p(specific)"
$ws.Range("B15").Value = 0.9993281923871153
$ws.Range("C15").Value = 0.9985745993077032
$ws.Range("D15").Value = 0.9995282850021244
$ws.Range("E15").Value = ""

# ---- sheet: lamda1 ----
$ws = $wb.Worksheets.Item("lamda1")

# Insert two blank rows before row 5 (will become "pi" and "Keq")
$ws.Range("A5:E6").Insert()

# Remove the old trailing "pi" row, now shifted down to row 16
$ws.Rows.Item(16).Delete()

# Give the two new label cells (A5:A6) the same bold/bordered format as the
# other row labels in column A (copy format from A4, the "lamda" label cell)
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

$ws.Range("A2").Value = "gain"
$ws.Range("B2").Value = 7.091339770682374
$ws.Range("C2").Value = 7.070818935122439
$ws.Range("D2").Value = 7.113135877405636
$ws.Range("E2").Value = 7

$ws.Range("A3").Value = "proximity"
$ws.Range("B3").Value = 0.2171760650968842
$ws.Range("C3").Value = 0.2053141933922834
$ws.Range("D3").Value = 0.2278114621634793
$ws.Range("E3").Value = 0.2

$ws.Range("A4").Value = "lamda"
$ws.Range("B4").Value = 0.9913055644816908
$ws.Range("C4").Value = 0.9630704462345296
$ws.Range("D4").Value = 1.018960305113115
$ws.Range("E4").Value = 1

$ws.Range("A5").Value = "pi"
$ws.Range("B5").Value = 0.1737028249976817
$ws.Range("C5").Value = 0.1592361811506343
$ws.Range("D5").Value = 0.1900410159483033
$ws.Range("E5").Value = 0.15

$ws.Range("A6").Value = "Keq"
$ws.Range("B6").Value = 0.2113543413197227
$ws.Range("C6").Value = 0.1893946654237603
$ws.Range("D6").Value = 0.2346304610812667
$ws.Range("E6").Value = ""

$ws.Range("A7").Value = "SNR"
$ws.Range("B7").Value = 3.948417826326017
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 3.75889905491011

$ws.Range("A8").Value = "MCC"
$ws.Range("B8").Value = 0.9622567739333436
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""

$ws.Range("A9").Value = "Recall"
$ws.Range("B9").Value = 0.9904306220095693
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""

$ws.Range("A10").Value = "Precision"
$ws.Range("B10").Value = 0.9473684210526315
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""

$ws.Range("A11").Value = "TN"
$ws.Range("B11").Value = 2059
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""

$ws.Range("A12").Value = "FP"
$ws.Range("B12").Value = 23
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""

$ws.Range("A13").Value = "FN"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""

$ws.Range("A14").Value = "TP"
$ws.Range("B14").Value = 414
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""

$ws.Range("A15").Value = "p(specific)"
$ws.Range("B15").Value = 0.9885547964978505
$ws.Range("C15").Value = 0.975139958193761
$ws.Range("D15").Value = 0.9940717758244612
$ws.Range("E15").Value = ""

# ---- sheet: lamda0.01 ----
$ws = $wb.Worksheets.Item("lamda0.01")

# Insert two blank rows before row 5 (will become "pi" and "Keq")
$ws.Range("A5:E6").Insert()

# Remove the old trailing "pi" row, now shifted down to row 16
$ws.Rows.Item(16).Delete()

# Give the two new label cells (A5:A6) the same bold/bordered format as the
# other row labels in column A (copy format from A4, the "lamda" label cell)
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

$ws.Range("A2").Value = "gain"
$ws.Range("B2").Value = 7.006089850771027
$ws.Range("C2").Value = 6.985053242317973
$ws.Range("D2").Value = 7.028435075021673
$ws.Range("E2").Value = 7

$ws.Range("A3").Value = "proximity"
$ws.Range("B3").Value = 0.2027042158313969
$ws.Range("C3").Value = 0.192454572273346
$ws.Range("D3").Value = 0.2119043999479361
$ws.Range("E3").Value = 0.2

$ws.Range("A4").Value = "lamda"
$ws.Range("B4").Value = 0.01039546471149803
$ws.Range("C4").Value = 0.007691950981008
$ws.Range("D4").Value = 0.0133870759672566
$ws.Range("E4").Value = 0.01

$ws.Range("A5").Value = "pi"
$ws.Range("B5").Value = 0.1732564768707913
$ws.Range("C5").Value = 0.1584778611899527
$ws.Range("D5").Value = 0.1899691401107213
$ws.Range("E5").Value = 0.15

$ws.Range("A6").Value = "Keq"
$ws.Range("B6").Value = 0.2107273659127262
$ws.Range("C6").Value = 0.1883228670989748
$ws.Range("D6").Value = 0.234520911207085
$ws.Range("E6").Value = ""

$ws.Range("A7").Value = "SNR"
$ws.Range("B7").Value = 3.749766130642911
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 3.75889905491011

$ws.Range("A8").Value = "MCC"
$ws.Range("B8").Value = 0.9985660739664485
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""

$ws.Range("A9").Value = "Recall"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""

$ws.Range("A10").Value = "Precision"
$ws.Range("B10").Value = 0.9976133651551312
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""

$ws.Range("A11").Value = "TN"
$ws.Range("B11").Value = 2081
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""

$ws.Range("A12").Value = "FP"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""

$ws.Range("A13").Value = "FN"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""

$ws.Range("A14").Value = "TP"
$ws.Range("B14").Value = 418
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""

$ws.Range("A15").Value = "p(specific)"
$ws.Range("B15").Value = 0.9998729624295279
$ws.Range("C15").Value = 0.9997184600760852
$ws.Range("D15").Value = 0.9999131081282366
$ws.Range("E15").Value = ""

# ---- sheet: lamda0.5 ----
$ws = $wb.Worksheets.Item("lamda0.5")

# Insert two blank rows before row 5 (will become "pi" and "Keq")
$ws.Range("A5:E6").Insert()

# Remove the old trailing "pi" row, now shifted down to row 16
$ws.Rows.Item(16).Delete()

# Give the two new label cells (A5:A6) the same bold/bordered format as the
# other row labels in column A (copy format from A4, the "lamda" label cell)
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

$ws.Range("A2").Value = "gain"
$ws.Range("B2").Value = 7.093184567093628
$ws.Range("C2").Value = 7.07258619372312
$ws.Range("D2").Value = 7.115063148529826
$ws.Range("E2").Value = 7

$ws.Range("A3").Value = "proximity"
$ws.Range("B3").Value = 0.2152565921174508
$ws.Range("C3").Value = 0.2032320271264699
$ws.Range("D3").Value = 0.2260526857329364
$ws.Range("E3").Value = 0.2

$ws.Range("A4").Value = "lamda"
$ws.Range("B4").Value = 0.4914945559425632
$ws.Range("C4").Value = 0.4733921863139085
$ws.Range("D4").Value = 0.5092979513473626
$ws.Range("E4").Value = 0.5

$ws.Range("A5").Value = "pi"
$ws.Range("B5").Value = 0.1732984386545297
$ws.Range("C5").Value = 0.1587799969581947
$ws.Range("D5").Value = 0.1897036573362159
$ws.Range("E5").Value = 0.15

$ws.Range("A6").Value = "Keq"
$ws.Range("B6").Value = 0.210766110079728
$ws.Range("C6").Value = 0.188749669744199
$ws.Range("D6").Value = 0.2341164349915773
$ws.Range("E6").Value = ""

$ws.Range("A7").Value = "SNR"
$ws.Range("B7").Value = 3.843139606860399
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 3.75889905491011

$ws.Range("A8").Value = "MCC"
$ws.Range("B8").Value = 0.971955289238697
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""

$ws.Range("A9").Value = "Recall"
$ws.Range("B9").Value = 0.9952153110047847
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""

$ws.Range("A10").Value = "Precision"
$ws.Range("B10").Value = 0.9585253456221198
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""

$ws.Range("A11").Value = "TN"
$ws.Range("B11").Value = 2064
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""

$ws.Range("A12").Value = "FP"
$ws.Range("B12").Value = 18
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""

$ws.Range("A13").Value = "FN"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""

$ws.Range("A14").Value = "TP"
$ws.Range("B14").Value = 416
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""

$ws.Range("A15").Value = "p(specific)"
$ws.Range("B15").Value = 0.9940361218491627
$ws.Range("C15").Value = 0.9870758688439624
$ws.Range("D15").Value = 0.997106970967898
$ws.Range("E15").Value = ""

# ---- sheet: lamda0.15 ----
$ws = $wb.Worksheets.Item("lamda0.15")

# Insert two blank rows before row 5 (will become "pi" and "Keq")
$ws.Range("A5:E6").Insert()

# Remove the old trailing "pi" row, now shifted down to row 16
$ws.Rows.Item(16).Delete()

# Give the two new label cells (A5:A6) the same bold/bordered format as the
# other row labels in column A (copy format from A4, the "lamda" label cell)
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

$ws.Range("A2").Value = "gain"
$ws.Range("B2").Value = 7.002403255994217
$ws.Range("C2").Value = 6.981462514546145
$ws.Range("D2").Value = 7.02464650658167
$ws.Range("E2").Value = 7

$ws.Range("A3").Value = "proximity"
$ws.Range("B3").Value = 0.2082011137697973
$ws.Range("C3").Value = 0.1980091052580437
$ws.Range("D3").Value = 0.217351398503064
$ws.Range("E3").Value = 0.2

$ws.Range("A4").Value = "lamda"
$ws.Range("B4").Value = 0.1542698026731528
$ws.Range("C4").Value = 0.1433309674317986
$ws.Range("D4").Value = 0.1652123308767115
$ws.Range("E4").Value = 0.15

$ws.Range("A5").Value = "pi"
$ws.Range("B5").Value = 0.172155532464627
$ws.Range("C5").Value = 0.1584130878392987
$ws.Range("D5").Value = 0.1876643266692657
$ws.Range("E5").Value = 0.15

$ws.Range("A6").Value = "Keq"
$ws.Range("B6").Value = 0.2090268040384629
$ws.Range("C6").Value = 0.1882314066731731
$ws.Range("D6").Value = 0.2310182362033923
$ws.Range("E6").Value = ""

$ws.Range("A7").Value = "SNR"
$ws.Range("B7").Value = 3.787464886780817
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 3.75889905491011

$ws.Range("A8").Value = "MCC"
$ws.Range("B8").Value = 0.9872644337931913
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""

$ws.Range("A9").Value = "Recall"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""

$ws.Range("A10").Value = "Precision"
$ws.Range("B10").Value = 0.9789227166276346
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""

$ws.Range("A11").Value = "TN"
$ws.Range("B11").Value = 2073
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""

$ws.Range("A12").Value = "FP"
$ws.Range("B12").Value = 9
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""

$ws.Range("A13").Value = "FN"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""

$ws.Range("A14").Value = "TP"
$ws.Range("B14").Value = 418
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""

$ws.Range("A15").Value = "p(specific)"
$ws.Range("B15").Value = 0.9980787137566742
$ws.Range("C15").Value = 0.9959616337693512
$ws.Range("D15").Value = 0.9986518429881251
$ws.Range("E15").Value = ""

$excel.CutCopyMode = 0
